$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("Across","Abuse","Actor","Acute","Admit","Adopt","Adult","Africa","After","Almond","America")

for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $values[$i]
}
